$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New shared text used to replace "User is on the MakeMyTrip home page" in E2 and E3
$newText = "User needs to stay login and is on the MakeMyTrip home page"

# Update E2 and E3 values
$ws.Range("E2").Value = $newText
$ws.Range("E3").Value = $newText

# E3 should end up wrapped (matches E2's existing style)
$ws.Range("E3").WrapText = $true

# F5 picks up the Times New Roman font (like E5) but keeps vertical-center + wrap,
# with general (not centered) horizontal alignment - copy E5's format then adjust wrap
$ws.Range("E5").Copy()
$ws.Range("F5").PasteSpecial(-4122)
$ws.Range("F5").HorizontalAlignment = 1
$ws.Range("F5").WrapText = $true

# Update sheet view: scroll position and selection
$ws.Application.ActiveWindow.ScrollColumn = 4
$ws.Range("E3").Select()
